# Add a new wishlist entry as the next row after the last used row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Arráncame la vida"
$ws.Range("B16").Value = "Angeles Mastretta"

# Column C is left blank for this entry, same as the existing rows above it
# (e.g. C8:C15) which are empty text cells rather than fully absent cells.
# Assigning a lone apostrophe forces a text-typed empty cell (Excel's
# "quote prefix" trick), then resetting the style back to Normal clears the
# quote-prefix formatting flag that the apostrophe entry would otherwise add.
$ws.Range("C16").Value = "'"
$ws.Range("C16").Style = "Normal"
